{"js": "// Replace the tail of the \".env\" instructions paragraph and drop the\n// paragraphs that used to list the literal environment-variable values.\n//\n// Before:\n//   \"...Cree un archivo .env en la carpeta backend/ con el siguiente contenido:\"\n//   \"CORS_ALLOWED_ORIGINS=...\"\n//   \"DATABASE_URL=...\"\n//   \"DEBUG=False\"\n//   \"DEFAULT_FROM_EMAIL=...\"\n//   \"DJANGO_ALLOWED_HOSTS=...\"\n//   \"EMAIL_HOST_USER=apikey\"\n//   \"FRONTEND_URL=...\"\n//   \"SECRET_KEY=...\"\n//   \"SENDGRID_API_KEY=...\"\n// After:\n//   \"...Cree un archivo .env en la carpeta backend/ y pide a los programadores el contenido.\"\n\nconst body = context.document.body;\n\n// 1) Swap the trailing sentence of the intro paragraph in place so the run's\n//    formatting (rPr) is preserved.\nconst tail = body.search(\" con el siguiente contenido:\", { matchCase: true });\ntail.load(\"items\");\nawait context.sync();\n\nif (tail.items.length === 0) {\n  throw new Error(\"Could not find the ' con el siguiente contenido:' text to replace.\");\n}\ntail.items[0].insertText(\" y pide a los programadores el contenido.\", \"Replace\");\nawait context.sync();\n\n// 2) Remove the now-obsolete paragraphs that used to dump the raw .env\n//    contents (from the CORS_ALLOWED_ORIGINS line through the\n//    SENDGRID_API_KEY line).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst startMarker = \"CORS_ALLOWED_ORIGINS=\";\nconst endMarker = \"SENDGRID_API_KEY=\";\n\nlet startIndex = -1;\nlet endIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (startIndex === -1 && text.indexOf(startMarker) === 0) {\n    startIndex = i;\n  }\n  if (text.indexOf(endMarker) === 0) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1 || endIndex < startIndex) {\n  throw new Error(\"Could not locate the block of .env value paragraphs to delete.\");\n}\n\nfor (let i = endIndex; i >= startIndex; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Fix Manual de Instalacion.docx\n#\n# The \".env\" setup paragraph used to end with \"con el siguiente contenido:\"\n# and was followed by paragraphs dumping the literal (secret-laden) .env\n# values. Replace the closing sentence and drop the leaked values instead.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the trailing sentence of the intro paragraph, keeping the rest\n#    of that paragraph (and its run formatting) untouched.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" con el siguiente contenido:\"\n$find.Replacement.Text = \" y pide a los programadores el contenido.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Remove the now-obsolete paragraphs that used to list the raw .env\n#    contents, from the CORS_ALLOWED_ORIGINS line through the\n#    SENDGRID_API_KEY line (inclusive).\n$startPara = $null\n$endPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($startPara -eq $null -and $t.StartsWith(\"CORS_ALLOWED_ORIGINS=\")) {\n        $startPara = $i\n    }\n    if ($t.StartsWith(\"SENDGRID_API_KEY=\")) {\n        $endPara = $i\n        break\n    }\n}\n\nif ($startPara -ne $null -and $endPara -ne $null) {\n    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start\n    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End\n    $deleteRange = $d.Range($rangeStart, $rangeEnd)\n    $deleteRange.Delete()\n}\n"}
